# daily auto push: 2025-10-10 13:36 UTC
# Append the new daily row (row 90) to the sheet, mirroring the existing
# data rows: date/weekday stored as plain text (not auto-converted to a
# date serial), time-of-day and ranking stored as numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 90

# Column A holds "2025/10/10" as literal text (matching every other date
# cell in the column). Force text storage via NumberFormat "@" before the
# assignment so Excel doesn't auto-detect it as a date serial, then reset
# the style back to Normal so no stray per-cell formatting is left behind.
$cellA = $ws.Cells.Item($row, 1)
$cellA.NumberFormat = "@"
$cellA.Value = "2025/10/10"
$cellA.Style = "Normal"

$ws.Cells.Item($row, 2).Value = "金"
$ws.Cells.Item($row, 3).Value = 20
$ws.Cells.Item($row, 4).Value = 34
